$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the three previously-missing "PDF do Projeto" filenames ---
# (rows identified by their "Projeto recebido (Autor)" value, before the table is re-sorted)

# Row where D = "Ailton Serrão Esquerdo"
$ws.Range("F3").Value = "Ailton_Serrão_Esquerdo.pdf"
$ws.Range("E3").Copy()
$ws.Range("F3").PasteSpecial(-4122)

# Row where D = "Alessandra Freixo Braga"
$ws.Range("F4").Value = "Alessandra_Freixo_Braga.pdf"
$ws.Range("E4").Copy()
$ws.Range("F4").PasteSpecial(-4122)

# Row where D = "Haroldo Andre Bastos da Silva"
$ws.Range("F17").Value = "Harodo_André_bastos_da_Silva.pdf"
$ws.Range("E17").Copy()
$ws.Range("F17").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Re-sort the whole table (A2:F35) by column A ascending instead of column D ---
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A35"))
$ws.Sort.SetRange($ws.Range("A1:F35"))
$ws.Sort.Header = 1
$ws.Sort.Apply()

# --- Move the active selection ---
$ws.Range("D13").Select()
